$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert 3 new columns before column B (shifts everything right by 3,
#    matching the diff's column-shift pattern across cols/row cells/merges).
$ws.Columns("B:D").Insert()

# 2. Header row (row 1): give B1:D1 the same look as A1 (Calibri font,
#    centered, wrapped) before filling in the new header text/rich runs.
$ws.Range("A1").Copy()
$ws.Range("B1:D1").PasteSpecial(-4122)

# day1 + "总分" (second run in 宋体)
$ws.Range("B1").Value = "day1总分"
$ws.Range("B1").Characters(5,2).Font.Name = "宋体"

# "git" + "组员互相平分" + " <newline>" + "第一次"
$ws.Range("C1").Value = "git组员互相平分 `n第一次"
$ws.Range("C1").Characters(4,6).Font.Name = "宋体"
$ws.Range("C1").Characters(10,2).Font.Name = "Calibri"
$ws.Range("C1").Characters(12,3).Font.Name = "宋体"

# "git " + "记分组评分<newline>第一次"
$ws.Range("D1").Value = "git 记分组评分`n第一次"
$ws.Range("D1").Characters(5,9).Font.Name = "宋体"

# 3. Data rows (2:36): blank cells styled like the rest of the score grid
#    (centered, same font as column A's roster entries).
$ws.Range("A14").Copy()
$ws.Range("B2:D36").PasteSpecial(-4122)
$ws.Range("B2:D36").ClearContents()

# 4. Column widths for the 3 new columns, matching column A's width as
#    closely as this engine's character-width quantization allows.
$ws.Columns("B:D").ColumnWidth = 17.75

# 5. Selection ends on B2, as in the saved file.
$ws.Range("B2").Select()
